$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "BB" column (date header + matching/forecast-revision values), rows 1-82
$bbValues = @{
    1 = 45986
    2 = 0.04899821040007168
    3 = 1.885178963001849
    4 = -1.307280175329765
    5 = 2.75267136845396
    6 = -4.471592960313714
    7 = -4.942495347902479
    8 = -5.434114574907241
    9 = 4.370618442157621
    10 = -0.9953340506219917
    11 = 5.692238679293155
    12 = 5.484876523251003
    13 = 0.1757137213762547
    14 = 1.626538719806248
    15 = 2.437334396728659
    16 = 1.566144859651857
    17 = 1.556824096124856
    18 = 0.21169683169569
    19 = 0.1582599013804469
    20 = -0.3253450194449812
    21 = -0.15853729477206
    22 = 0.2054976683197225
    23 = 0.9595582875050894
    24 = 1.503966953576466
    25 = 1.799837015295822
    26 = 0.6518403676065248
    27 = 0.6633823054011998
    28 = 0.923139910942723
    29 = 0.3869820931359413
    30 = 2.093916965767463
    31 = 1.580888475204972
    32 = 0.3008043112709089
    33 = 1.199598313222268
    34 = 0.4774400648527148
    35 = 1.507463254996111
    36 = 0.2179116434425623
    37 = 0.6120689161334525
    38 = 2.489390679284554
    39 = 0.5389418434166515
    40 = 2.40118094791471
    41 = 1.043009620608657
    42 = 1.830928398766659
    43 = -0.3674870133197601
    44 = 1.62717758729876
    45 = 1.117271732844245
    46 = 0.6525147083449099
    47 = 1.4
    48 = -0.3
    49 = -0.3
    50 = 0.1
    51 = -1.138880770453937
    52 = -16.88491062648744
    53 = 9.224715108933083
    54 = 3.283355339827622
    55 = 4.432584407022276
    56 = 2.509693347214139
    57 = -0.4381048169788073
    58 = 4.106981763725997
    59 = 0.3842995656585515
    60 = 2.277966437795897
    61 = 2.507553358214992
    62 = -2.399190900254823
    63 = -1.115644072253531
    64 = 0.1689348086957096
    65 = -1.814969742946232
    66 = -1.586779238813989
    67 = -0.2996177924633514
    68 = 2.757652919539751
    69 = -0.08601690538415596
    70 = -0.6801011570971838
    71 = 1.538981993999982
    72 = 1.68501852020853
    73 = 0.03331000006224372
    74 = 0.719473196736553
    75 = 0.6752544302532623
    76 = 0.6532787989117363
    77 = 0.6588538323071631
    78 = 0.6588609317435923
    79 = 0.6586210043061691
    80 = 0.6586608941090445
    81 = 0.6586644859197338
    82 = 0.6586621760120679
}
foreach ($row in $bbValues.Keys) {
    $ws.Cells.Item($row, 54).Value = $bbValues[$row]
}

# New row 83 (next forecast quarter)
$ws.Range("A83").Value = 46934
$ws.Range("BB83").Value = 0.6586624091680092

# BB1 (a header date) and A83 (a date in column A) need the same date
# number-format/border/bold style ("s=1") used by the rest of the date cells,
# so copy it over from an existing, already-styled date cell.
$ws.Range("BA1").Copy()
$ws.Range("BB1").PasteSpecial(-4122)
$ws.Range("A82").Copy()
$ws.Range("A83").PasteSpecial(-4122)
$excel.CutCopyMode = 0
